$wb = $excel.ActiveWorkbook

# 1) Rename header cells on existing sheets
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 2) Add the new "PO Forecast" sheet at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "PO Forecast"

# 3) Write header row
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Copy header formatting (bold, centered, bordered) from an existing header cell
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Copy date-column formatting from an existing date cell
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A75").PasteSpecial(-4122)  # xlPasteFormats

# 4) Populate the forecast data rows (rows 2-75)
$data = New-Object 'object[,]' 74,4
$data[0,0] = 44934.99999999999
$data[0,1] = 52
$data[0,2] = -20.99848901921702
$data[0,3] = 121.5955254002123
$data[1,0] = 44941.99999999999
$data[1,1] = 53
$data[1,2] = -21.25030725646484
$data[1,3] = 126.7630671880579
$data[2,0] = 44962.99999999999
$data[2,1] = 54
$data[2,2] = -19.21252364849919
$data[2,3] = 125.2755783868481
$data[3,0] = 44990.99999999999
$data[3,1] = 56
$data[3,2] = -15.48969759150752
$data[3,3] = 133.7283970731424
$data[4,0] = 44997.99999999999
$data[4,1] = 56
$data[4,2] = -15.1355135170434
$data[4,3] = 133.7943857698895
$data[5,0] = 45004.99999999999
$data[5,1] = 57
$data[5,2] = -16.49892224961599
$data[5,3] = 127.1889818596016
$data[6,0] = 45011.99999999999
$data[6,1] = 57
$data[6,2] = -11.34794395678288
$data[6,3] = 132.7599790517683
$data[7,0] = 45018.99999999999
$data[7,1] = 58
$data[7,2] = -5.391442197781912
$data[7,3] = 128.9268748464916
$data[8,0] = 45025.99999999999
$data[8,1] = 58
$data[8,2] = -16.3480195445139
$data[8,3] = 130.3488242477637
$data[9,0] = 45032.99999999999
$data[9,1] = 59
$data[9,2] = -14.13908794505816
$data[9,3] = 128.7785441747996
$data[10,0] = 45039.99999999999
$data[10,1] = 59
$data[10,2] = -11.35706018765515
$data[10,3] = 130.7394198656195
$data[11,0] = 45046.99999999999
$data[11,1] = 60
$data[11,2] = -12.08059213108936
$data[11,3] = 133.5954265124785
$data[12,0] = 45060.99999999999
$data[12,1] = 61
$data[12,2] = -9.566498308665377
$data[12,3] = 135.5852685630874
$data[13,0] = 45067.99999999999
$data[13,1] = 61
$data[13,2] = -8.747925687957899
$data[13,3] = 137.4411620290095
$data[14,0] = 45081.99999999999
$data[14,1] = 62
$data[14,2] = -1.061560785003743
$data[14,3] = 141.1592588562246
$data[15,0] = 45088.99999999999
$data[15,1] = 63
$data[15,2] = -10.54996126515564
$data[15,3] = 128.4414001755519
$data[16,0] = 45102.99999999999
$data[16,1] = 64
$data[16,2] = -16.04046136605115
$data[16,3] = 129.1481027665276
$data[17,0] = 45109.99999999999
$data[17,1] = 64
$data[17,2] = -1.366052131173523
$data[17,3] = 139.511715995806
$data[18,0] = 45116.99999999999
$data[18,1] = 65
$data[18,2] = -12.9118267909601
$data[18,3] = 134.4968443821462
$data[19,0] = 45130.99999999999
$data[19,1] = 66
$data[19,2] = -8.514313927982837
$data[19,3] = 136.1686661468753
$data[20,0] = 45137.99999999999
$data[20,1] = 66
$data[20,2] = -7.870474639676005
$data[20,3] = 140.3417571090607
$data[21,0] = 45144.99999999999
$data[21,1] = 67
$data[21,2] = -12.07635906352437
$data[21,3] = 138.7011363260752
$data[22,0] = 45151.99999999999
$data[22,1] = 67
$data[22,2] = -4.958295725175137
$data[22,3] = 137.5973509691134
$data[23,0] = 45165.99999999999
$data[23,1] = 68
$data[23,2] = -4.619473033767568
$data[23,3] = 142.7469447388963
$data[24,0] = 45172.99999999999
$data[24,1] = 69
$data[24,2] = -8.705033211206905
$data[24,3] = 139.8395926036491
$data[25,0] = 45179.99999999999
$data[25,1] = 69
$data[25,2] = -4.337028293454381
$data[25,3] = 140.3171041722277
$data[26,0] = 45186.99999999999
$data[26,1] = 70
$data[26,2] = -1.484168767002346
$data[26,3] = 139.2244685854036
$data[27,0] = 45200.99999999999
$data[27,1] = 71
$data[27,2] = -0.08538266503827646
$data[27,3] = 140.1103942990447
$data[28,0] = 45207.99999999999
$data[28,1] = 71
$data[28,2] = -1.085124972742556
$data[28,3] = 141.6621941856974
$data[29,0] = 45214.99999999999
$data[29,1] = 72
$data[29,2] = -5.509939421693113
$data[29,3] = 147.1963585037458
$data[30,0] = 45221.99999999999
$data[30,1] = 72
$data[30,2] = -2.80944167420405
$data[30,3] = 142.2569452327969
$data[31,0] = 45228.99999999999
$data[31,1] = 73
$data[31,2] = -0.2280694937458957
$data[31,3] = 147.4847333642974
$data[32,0] = 45242.99999999999
$data[32,1] = 74
$data[32,2] = 2.676481847369051
$data[32,3] = 149.7116848802706
$data[33,0] = 45249.99999999999
$data[33,1] = 74
$data[33,2] = 4.438074031366477
$data[33,3] = 142.2887908774609
$data[34,0] = 45256.99999999999
$data[34,1] = 75
$data[34,2] = 2.455137648832554
$data[34,3] = 151.6652772509728
$data[35,0] = 45263.99999999999
$data[35,1] = 75
$data[35,2] = 7.251699229675038
$data[35,3] = 151.7126972364254
$data[36,0] = 45277.99999999999
$data[36,1] = 76
$data[36,2] = 1.467202022707045
$data[36,3] = 141.3335508995345
$data[37,0] = 45298.99999999999
$data[37,1] = 77
$data[37,2] = 8.794344528569859
$data[37,3] = 147.0856881594674
$data[38,0] = 45305.99999999999
$data[38,1] = 78
$data[38,2] = 8.573457193206369
$data[38,3] = 150.6755539010194
$data[39,0] = 45312.99999999999
$data[39,1] = 78
$data[39,2] = 6.062526072946838
$data[39,3] = 149.9916637894486
$data[40,0] = 45319.99999999999
$data[40,1] = 79
$data[40,2] = 11.84270655259256
$data[40,3] = 148.5972018980264
$data[41,0] = 45326.99999999999
$data[41,1] = 79
$data[41,2] = 7.412876836223893
$data[41,3] = 154.9231107885403
$data[42,0] = 45333.99999999999
$data[42,1] = 80
$data[42,2] = 5.967766533493407
$data[42,3] = 151.3017700317453
$data[43,0] = 45347.99999999999
$data[43,1] = 81
$data[43,2] = 10.07701484027931
$data[43,3] = 157.0899320966641
$data[44,0] = 45354.99999999999
$data[44,1] = 81
$data[44,2] = 19.41758487795244
$data[44,3] = 159.8405946055595
$data[45,0] = 45361.99999999999
$data[45,1] = 82
$data[45,2] = 12.64503171366688
$data[45,3] = 154.3221762262638
$data[46,0] = 45368.99999999999
$data[46,1] = 82
$data[46,2] = 5.469864704602981
$data[46,3] = 152.6786851634922
$data[47,0] = 45382.99999999999
$data[47,1] = 83
$data[47,2] = 13.65271599359811
$data[47,3] = 156.4894707915205
$data[48,0] = 45389.99999999999
$data[48,1] = 84
$data[48,2] = 14.04317362639197
$data[48,3] = 153.0952532833041
$data[49,0] = 45396.99999999999
$data[49,1] = 84
$data[49,2] = 15.01587430148759
$data[49,3] = 155.2448088029642
$data[50,0] = 45403.99999999999
$data[50,1] = 85
$data[50,2] = 21.92019495270311
$data[50,3] = 156.9663384244014
$data[51,0] = 45410.99999999999
$data[51,1] = 85
$data[51,2] = 15.15663538372405
$data[51,3] = 159.7915843486741
$data[52,0] = 45417.99999999999
$data[52,1] = 86
$data[52,2] = 10.7977203989782
$data[52,3] = 159.412831960436
$data[53,0] = 45424.99999999999
$data[53,1] = 86
$data[53,2] = 17.44086884762388
$data[53,3] = 163.1162718484794
$data[54,0] = 45431.99999999999
$data[54,1] = 87
$data[54,2] = 14.01557505755938
$data[54,3] = 156.6085415762335
$data[55,0] = 45438.99999999999
$data[55,1] = 87
$data[55,2] = 13.6474052945418
$data[55,3] = 163.915124980718
$data[56,0] = 45452.99999999999
$data[56,1] = 88
$data[56,2] = 21.37554946990387
$data[56,3] = 159.2773864408265
$data[57,0] = 45459.99999999999
$data[57,1] = 89
$data[57,2] = 17.22237743230095
$data[57,3] = 160.3872137282288
$data[58,0] = 45466.99999999999
$data[58,1] = 89
$data[58,2] = 13.37515258311059
$data[58,3] = 162.4117267397455
$data[59,0] = 45473.99999999999
$data[59,1] = 90
$data[59,2] = 18.08788812665768
$data[59,3] = 158.7214935500438
$data[60,0] = 45480.99999999999
$data[60,1] = 90
$data[60,2] = 15.79072277458045
$data[60,3] = 160.7363482395562
$data[61,0] = 45487.99999999999
$data[61,1] = 91
$data[61,2] = 18.03371148589585
$data[61,3] = 168.4038076741377
$data[62,0] = 45494.99999999999
$data[62,1] = 91
$data[62,2] = 21.63806880961916
$data[62,3] = 166.6603973099721
$data[63,0] = 45522.99999999999
$data[63,1] = 93
$data[63,2] = 17.60066236081373
$data[63,3] = 167.7642319517199
$data[64,0] = 45529.99999999999
$data[64,1] = 94
$data[64,2] = 17.82814860246367
$data[64,3] = 165.2367532751348
$data[65,0] = 45564.99999999999
$data[65,1] = 96
$data[65,2] = 26.91905053443889
$data[65,3] = 169.8883003047069
$data[66,0] = 45571.99999999999
$data[66,1] = 97
$data[66,2] = 28.42096875143348
$data[66,3] = 171.3323684043543
$data[67,0] = 45578.99999999999
$data[67,1] = 97
$data[67,2] = 23.00541147544425
$data[67,3] = 173.5968074135769
$data[68,0] = 45585.99999999999
$data[68,1] = 97
$data[68,2] = 26.77746722112089
$data[68,3] = 168.8904306308976
$data[69,0] = 45592.99999999999
$data[69,1] = 98
$data[69,2] = 23.84029367308956
$data[69,3] = 171.6489118482586
$data[70,0] = 45599.99999999999
$data[70,1] = 98
$data[70,2] = 28.34753875805991
$data[70,3] = 171.7706272074424
$data[71,0] = 45606.99999999999
$data[71,1] = 99
$data[71,2] = 27.52557043921307
$data[71,3] = 171.0143125585818
$data[72,0] = 45613.99999999999
$data[72,1] = 99
$data[72,2] = 28.76978043840086
$data[72,3] = 174.8399655734297
$data[73,0] = 45620.99999999999
$data[73,1] = 100
$data[73,2] = 31.43841199105998
$data[73,3] = 175.6042087793401

$newSheet.Range("A2:D75").Value = $data

Write-Host "PO Forecast sheet created and headers updated."
